# Rename the inline Pearson/BTEC logo pictures embedded in the section's
# headers and footers:
#   - footer 1 (wdHeaderFooterPrimary)   : image2.png -> image1.png
#   - footer 2 (wdHeaderFooterFirstPage) : image2.png -> image1.png
#   - header 2 (wdHeaderFooterFirstPage) : image1.jpg -> image2.jpg

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    $footer1 = $sec.Footers.Item(1)
    if ($footer1.Range.InlineShapes.Count -ge 1) {
        $footer1.Range.InlineShapes.Item(1).Name = "image1.png"
    }

    $footer2 = $sec.Footers.Item(2)
    if ($footer2.Range.InlineShapes.Count -ge 1) {
        $footer2.Range.InlineShapes.Item(1).Name = "image1.png"
    }

    $header2 = $sec.Headers.Item(2)
    if ($header2.Range.InlineShapes.Count -ge 1) {
        $header2.Range.InlineShapes.Item(1).Name = "image2.jpg"
    }
}
